# "version final sin errores"
# - Bump the CodeSystem "Version" property from 0.4.0 to 0.7.0
# - Remove the now-unused "Jurisdiction" / "Chile" property row entirely
#   (shifts every row below it up by one; the Concepts sheet is untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3: Property="Version", Value="0.4.0" -> "0.7.0")
$ws.Range("B3").Value = "0.7.0"

# Delete the entire "Jurisdiction" / "Chile" row (row 11)
$ws.Rows.Item(11).Delete()
